# Apply the recorded edit to the workbook:
#  1. Rename the "March" sheet to "mars".
#  2. Append a new row (row 6) to the expense table with:
#       Category=Clothing, Name=genser, Date=2023-03-20, Price=1000.0, Account=Checkings
#
# Dates / numeric-looking strings ("2023-03-20", "1000.0") must be written
# as literal text (matching the existing Date/Price columns, which are
# plain shared-string text, not real dates/numbers). To stop Excel's
# automatic type inference from turning them into a date serial / number we
# briefly mark the cell as Text ("@") before assigning the value, then
# clear the cell formatting again so the cell itself keeps the default
# (unstyled) look of every other cell in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet "March" -> "mars"
$ws.Name = "mars"

# 2. Add the new row of data (row 6)
$ws.Range("A6").Value = "Clothing"
$ws.Range("B6").Value = "genser"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "2023-03-20"
$ws.Range("C6").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1000.0"
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = "Checkings"
